$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.542896509170532
$ws.Range("B1").Value = 2.185683488845825
$ws.Range("C1").Value = -1
$ws.Range("D1").Value = 1.363738656044006
$ws.Range("E1").Value = 0.6506624221801758
